$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column width adjustments (closest achievable via COM pixel-quantized ColumnWidth)
$ws.Columns.Item(1).ColumnWidth = 15.666666666666666
$ws.Columns.Item(2).ColumnWidth = 14.833333333333334

# Cell value updates
$ws.Range("A1").Value = -0.14696094660155978
$ws.Range("B1").Value = 0.14690035879450392
$ws.Range("A2").Value = -0.12479585091835865
$ws.Range("B2").Value = 0.12454213909045375
$ws.Range("A3").Value = -0.063507141423812641
$ws.Range("B3").Value = 0.063218316575090583
$ws.Range("A4").Value = -0.055218316623916408
$ws.Range("B4").Value = 0.054947996580853697
$ws.Range("A5").Value = -0.051947996605709257
$ws.Range("B5").Value = 0.051046937150625027
$ws.Range("A6").Value = -0.024881430387829084
$ws.Range("B6").Value = 0.024650921759185707
$ws.Range("A7").Value = -0.0146509218251758
$ws.Range("B7").Value = 0.014606991015348569
$ws.Range("A8").Value = -0.0046069910822463811
$ws.Range("B8").Value = 0.0045573042672240049
$ws.Range("A9").Value = -0.0025573042935937984
$ws.Range("B9").Value = 0.0025249380246918918
$ws.Range("A10").Value = -0.00052493805139519623
$ws.Range("B10").Value = 0.00052476133849665985
$ws.Range("A11").Value = -0.02437902214033727
$ws.Range("B11").Value = 0.024353565187637649
$ws.Range("A12").Value = -0.020853565222766601
$ws.Range("B12").Value = 0.02066545821033472
$ws.Range("A13").Value = -0.017165458248634913
$ws.Range("B13").Value = 0.017079788591179756
$ws.Range("A14").Value = -0.0090797886534614847
$ws.Range("B14").Value = 0.0090517927421558397
$ws.Range("A15").Value = -0.0080517927689829349
$ws.Range("B15").Value = 0.0080337078951382423
$ws.Range("A16").Value = -0.006033707927612042
$ws.Range("B16").Value = 0.0060034356699847358
$ws.Range("A17").Value = -0.0040034357032023848
$ws.Range("B17").Value = 0.0039999999562976285
$ws.Range("A18").Value = -0.016104507907680699
$ws.Range("B18").Value = 0.016091433885488726
$ws.Range("A19").Value = -0.012091433906887161
$ws.Range("B19").Value = 0.012016770121363596
$ws.Range("A20").Value = -0.0080167701443532025
$ws.Range("B20").Value = 0.0080056944840674049
$ws.Range("A21").Value = -0.0040056945072821648
$ws.Range("B21").Value = 0.0039999999765898409
$ws.Range("A22").Value = -0.057034997693074629
$ws.Range("B22").Value = 0.056742156459280935
$ws.Range("A23").Value = -0.040491250574224757
$ws.Range("B23").Value = 0.040097485066238114
$ws.Range("A24").Value = -0.020097485180484931
$ws.Range("B24").Value = 0.01999999988424328
$ws.Range("A25").Value = -0.024165506787170443
$ws.Range("B25").Value = 0.0241444851581889
$ws.Range("A26").Value = -0.021644485185449369
$ws.Range("B26").Value = 0.021620520332094273
$ws.Range("A27").Value = -0.03403203307919167
$ws.Range("B27").Value = 0.033920064969370944
$ws.Range("A28").Value = -0.031920064997759567
$ws.Range("B28").Value = 0.031855015183148261
$ws.Range("A29").Value = -0.024855015239468869
$ws.Range("B29").Value = 0.024843988947177209
$ws.Range("A30").Value = 0.035156010723502717
$ws.Range("B30").Value = -0.035241731465617487
$ws.Range("A31").Value = 0.042241731411316152
$ws.Range("B31").Value = -0.042299004918849192
$ws.Range("A32").Value = -0.0040012889015912378
$ws.Range("B32").Value = 0.0039999999619944049
